# The source workbook (Y6_n.xlsx) had an accidental duplicate block of
# 24 rows (wavelengths 1485-1600 were written twice into the sheet, once
# at rows 225-248 and again at rows 249-272). This edit removes the
# duplicated second copy so the wavelength series becomes monotonic again
# (..., 1595, 1600, 1605, 1610, ...) and the used range shrinks from
# A1:B352 down to A1:B328.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 24 duplicated rows (former rows 249-272: wavelengths 1485-1600
# repeated). Everything below shifts up to close the gap.
$ws.Rows("249:272").Delete()

# Restore the on-screen cursor/scroll position to where the author left it
# after trimming the duplicate rows.
$ws.Range("J313").Select()
